$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (Word table row 1)
$t.Cell(1,1).Range.Text = "90÷5=18, 0"
$t.Cell(1,2).Range.Text = "50÷4=12, 2"
$t.Cell(1,3).Range.Text = "99÷4=24, 3"
$t.Cell(1,4).Range.Text = "68÷8=8, 4"
$t.Cell(1,5).Range.Text = "63÷4=15, 3"

# Row 2 (Word table row 5)
$t.Cell(5,1).Range.Text = "60÷5=12, 0"
$t.Cell(5,2).Range.Text = "58÷6=9, 4"
$t.Cell(5,3).Range.Text = "92÷6=15, 2"
$t.Cell(5,4).Range.Text = "91÷4=22, 3"
$t.Cell(5,5).Range.Text = "41÷3=13, 2"

# Row 3 (Word table row 9)
$t.Cell(9,1).Range.Text = "93÷3=31, 0"
$t.Cell(9,2).Range.Text = "31÷7=4, 3"
$t.Cell(9,3).Range.Text = "14÷3=4, 2"
$t.Cell(9,4).Range.Text = "63÷3=21, 0"
$t.Cell(9,5).Range.Text = "80÷8=10, 0"

# Row 4 (Word table row 13)
$t.Cell(13,1).Range.Text = "95÷3=31, 2"
$t.Cell(13,2).Range.Text = "54÷2=27, 0"
$t.Cell(13,3).Range.Text = "41÷4=10, 1"
$t.Cell(13,4).Range.Text = "43÷2=21, 1"
$t.Cell(13,5).Range.Text = "37÷7=5, 2"

# Row 5 (Word table row 17)
$t.Cell(17,1).Range.Text = "18÷4=4, 2"
$t.Cell(17,2).Range.Text = "35÷7=5, 0"
$t.Cell(17,3).Range.Text = "63÷6=10, 3"
$t.Cell(17,4).Range.Text = "73÷2=36, 1"
$t.Cell(17,5).Range.Text = "73÷2=36, 1"
